$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns are treated as text so numeric-looking
# strings (e.g. "589.68", "0.999") are not auto-converted to numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.407.20'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.502.42'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '589.68'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.18'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.73'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +7.30%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.387'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.096.57'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.35%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.32%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.503.02'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.328.06'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.33'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.77%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.03'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.78'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.28%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '387.20'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.641.51'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.40'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.29%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.75'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.38%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.84%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.992'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.34'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.82%  '
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.52'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.96%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.25'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.15'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.530.96'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.56%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.66%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +4.35%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.55'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '164.88'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0789'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.808'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.42'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.34%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.19'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '24.27'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -5.09%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.27%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.439.67'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.58%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.88%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.922'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.25%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.29%  '
